$wb = $excel.ActiveWorkbook

# Sheet: 展览
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F7").Value = 385
$ws.Range("F9").Value = 1184
$ws.Range("F10").Value = 702
$ws.Range("F11").Value = 466
$ws.Range("F12").Value = 202
$ws.Range("F13").Value = 767
$ws.Range("F14").Value = 85
$ws.Range("F16").Value = 16
$ws.Range("G16").NumberFormat = "@"
$ws.Range("G16").Value = "68"
$ws.Range("G16").Style = "Normal"
$ws.Range("F17").Value = 188
$ws.Range("F18").Value = 255
$ws.Range("F20").Value = 308
$ws.Range("F21").Value = 157
$ws.Range("F22").Value = 1569
$ws.Range("F23").Value = 140
$ws.Range("F24").Value = 75
$ws.Range("F27").Value = 2229
$ws.Range("F28").Value = 121
$ws.Range("F29").Value = 13
$ws.Range("F34").Value = 73

# Sheet: 演出
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F15").Value = 458
$ws.Range("F16").Value = 169

# Sheet: 本地生活
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F7").Value = 599

# Sheet: 全部类型
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F8").Value = 599
$ws.Range("F16").Value = 385
$ws.Range("F18").Value = 1184
$ws.Range("F19").Value = 702
$ws.Range("F20").Value = 466
$ws.Range("F23").Value = 767
$ws.Range("F24").Value = 85
$ws.Range("F29").Value = 188
$ws.Range("F30").Value = 255
$ws.Range("F32").Value = 308
$ws.Range("F34").Value = 157
$ws.Range("F35").Value = 1569
$ws.Range("F36").Value = 140
$ws.Range("F37").Value = 458
$ws.Range("F38").Value = 75
$ws.Range("F41").Value = 2229
$ws.Range("F42").Value = 169
$ws.Range("F43").Value = 121
$ws.Range("F48").Value = 73
